$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.390.10"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.664.64"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.23"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3964"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3910"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.14"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.399"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08597"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.41"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.315"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001357"
$ws.Range("E15").Value = "  +5.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.877"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.664.41"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.29"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06986"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.53"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.002"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.75"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.420.85"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.425"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.050"
$ws.Range("E26").Value = "  +11.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.53"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.67"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.438"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.106"
$ws.Range("E31").Value = "  -9.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.523"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.845.68"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.066"
$ws.Range("E34").Value = "  +8.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08274"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03026"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.900"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2768"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.05"
$ws.Range("E39").Value = "  +9.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09232"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7746"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.80"
$ws.Range("E42").Value = "  +4.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.441"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.54"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7115"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.538"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.135"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08452"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.32"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.270"
$ws.Range("E51").Value = "  +0.30%  "
